# Sunil Narine.xlsx — update runs/balls/fours/sixes (columns C:F) for rows 2-9.
# The stat rows got reshuffled between different matches; each destination
# row below is populated with the values that used to live on a different
# row of the same table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as text cells (matches the sheet's existing text-stored numbers)
# instead of letting Excel auto-convert numeric-looking strings to numbers.
$ws.Range("C2:F9").NumberFormat = "@"

$newValues = @{
    2 = @("3", "5", "0", "0")
    3 = @("0", "2", "0", "0")
    4 = @("6", "4", "1", "0")
    5 = @("17", "9", "1", "1")
    6 = @("9", "10", "0", "1")
    7 = @("7", "7", "0", "1")
    8 = @("0", "2", "0", "0")
    9 = @("64", "32", "6", "4")
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
}
